$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-7:
# serial date 45204 (2023-10-05) -> 45207 (2023-10-08)
foreach ($r in 2..7) {
    $ws.Cells.Item($r, 3).Value = 45207
}
